# Apply updates to the "F" column (view/attendance counts) across the
# workbook's sheets, matching the regenerated gh-pages data snapshot.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value  = 501
$ws1.Range("F14").Value = 1470
$ws1.Range("F19").Value = 227
$ws1.Range("F20").Value = 227
$ws1.Range("F24").Value = 5479
$ws1.Range("F25").Value = 4668
$ws1.Range("F31").Value = 1225
$ws1.Range("F37").Value = 208

# --- Sheet "本地生活" (Local life) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 2406

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F8").Value  = 501
$ws4.Range("F17").Value = 1470
$ws4.Range("F22").Value = 227
$ws4.Range("F23").Value = 227
$ws4.Range("F29").Value = 5479
$ws4.Range("F30").Value = 4668
$ws4.Range("F33").Value = 1225
$ws4.Range("F44").Value = 208
